$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before existing row 2 (shifts old data rows 2-21 down to 11-30)
$ws.Rows("2:10").Insert()
# The inserted rows pick up formatting from the row above (the bold header);
# clear that so the new data rows stay unformatted like the rest of the data.
$ws.Range("A2:C10").ClearFormats()

# New data for rows 2-10 (x, y, z)
$newTop = @(
    @(-2.896898627281189, 7.169353723526001, -0.4825034886598586),
    @(-3.050878047943115, 7.165829300880432, -0.3550609424710275),
    @(-3.056754767894745, 7.225212574005127, -0.4760921187698841),
    @(-3.037489891052246, 7.236634731292725, -0.4997432827949525),
    @(-2.916959762573243, 7.203977525234222, -0.616998553276062),
    @(-3.009585857391357, 7.19498348236084, -0.6860059350728989),
    @(-2.86443132162094, 7.121285438537598, -0.5276834592223163),
    @(-2.889585494995118, 7.118069887161255, -0.4351722449064255),
    @(-3.056696653366089, 7.102567493915558, -0.6014280728995802)
)

for ($i = 0; $i -lt $newTop.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTop[$i][2]
}

# Append a new row at the bottom (row 31)
$ws.Cells.Item(31, 1).Value = 2.491997003555297
$ws.Cells.Item(31, 2).Value = 6.562706351280213
$ws.Cells.Item(31, 3).Value = -1.225150167942047
